$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "GSTIN"
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "Contact Person"
$ws.Range("E1").Value = "Mobile Number"
$ws.Range("F1").Value = "Email ID"

# Format mobile-number column as text so numeric-looking values are kept as strings
$ws.Range("E2:E6").NumberFormat = "@"

# Row 2 - America Construction Ltd
$ws.Range("A2").Value = "America Construction Ltd"
$ws.Range("B2").Value = "29ABCDE1234F1Z5"
$ws.Range("C2").Value = "123 Main Street, Mumbai, Maharashtra 400001"
$ws.Range("D2").Value = "Rajesh Kumar"
$ws.Range("E2").Value = "9876543210"
$ws.Range("F2").Value = "rajesh.kumar@abcconstruction.com"

# Row 3 - Australia Builders Pvt Ltd
$ws.Range("A3").Value = "Australia Builders Pvt Ltd"
$ws.Range("B3").Value = "30FGHIJ5678K2L6"
$ws.Range("C3").Value = "456 Park Avenue, Delhi, Delhi 110001"
$ws.Range("D3").Value = "Priya Sharma"
$ws.Range("E3").Value = "9876543211"
$ws.Range("F3").Value = "priya.sharma@xyzbuilders.com"

# Row 4 - Canada Materials Co
$ws.Range("A4").Value = "Canada Materials Co"
$ws.Range("B4").Value = "27KLMNO9012P3Q7"
$ws.Range("C4").Value = "789 Industrial Area, Bangalore, Karnataka 560001"
$ws.Range("D4").Value = "Amit Patel"
$ws.Range("E4").Value = "9876543212"
$ws.Range("F4").Value = "amit.patel@defmaterials.com"

# Row 5 - India Suppliers (no GSTIN, same as before)
$ws.Range("A5").Value = "India Suppliers"
$ws.Range("C5").Value = "321 Trade Center, Pune, Maharashtra 411001"
$ws.Range("D5").Value = "Sneha Desai"
$ws.Range("E5").Value = "9876543213"
$ws.Range("F5").Value = "sneha.desai@ghisuppliers.com"

# Row 6 - UK Enterprises
$ws.Range("A6").Value = "UK Enterprises"
$ws.Range("B6").Value = "24RSTUV3456W4X8"
$ws.Range("C6").Value = "654 Business Park, Hyderabad, Telangana 500001"
$ws.Range("D6").Value = "Vikram Singh"
$ws.Range("E6").Value = "9876543214"
$ws.Range("F6").Value = "vikram.singh@jklenterprises.com"

# Reset column B formatting (removes its explicit custom width so it reverts
# to the sheet default, matching the dropped <col> definition for column B)
$ws.Columns.Item(2).ClearFormats()

# B5 was never populated (row 5 only has Name/Address, no GSTIN) - make sure no
# stray empty cell record is left behind after touching the column formatting
$ws.Range("B5").ClearContents()

# Move selection to A7
[void]$ws.Range("A7").Select()
